$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Campaigns")

$values = @(
  "Campaign Name",
  "apple»artichoke»apple»carrot»carrot»robin»hamster»robin»hamster»robin»thor»hamster»hamster",
  "grape»artichoke»grape»artichoke»artichoke»robin»hamster»robin»hamster»robin»thor»hamster»hamster",
  "apple»carrot»apple»carrot»carrot»robin»hamster»robin»antelope»robin»thor»hamster»hamster",
  "grape»artichoke»grape»artichoke»carrot»robin»hamster»robin»hamster»robin»thor»hamster»hamster",
  "grape»artichoke»grape»carrot»carrot»robin»hamster»robin»hamster»robin»thor»hamster»hamster",
  "apple»carrot»banana»carrot»carrot»hawk»hamster»hawk»hamster»robin»spider-man»antelope»antelope",
  "grape»artichoke»apple»carrot»carrot»robin»hamster»robin»hamster»robin»thor»hamster»hamster",
  "grape»carrot»apple»carrot»carrot»robin»hamster»robin»hamster»robin»thor»hamster»hamster",
  "banana»carrot»apple»carrot»carrot»robin»hamster»robin»hamster»robin»thor»hamster»hamster",
  "apple»carrot»apple»carrot»carrot»robin»hamster»robin»hamster»robin»thor»hamster»hamster",
  "CV4/CFH/683_Chevrolet Retail SC _Chevy Full Line HA_Austin Hispanic Regional Plan_BPR_SCR_GNM_2018",
  "1LM/PRO/683_Chevy North East LMA_PROVIDENCE / NEW BEDFORD_Q2 LMA_BCN_NER_GNM_2018",
  "apple»carrot»banana»kale»kale»robin»antelope»hawk»hamster»hawk»spider-man»antelope»antelope",
  "apple»carrot»banana»kale»kale»hawk»antelope»hawk»antelope»hawk»spider-man»antelope»antelope"
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $values[$i]
}
